# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "46.149.93"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "2.587.95"
$ws.Range("E3").Value = "  +8.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.38"
$ws.Range("E5").Value = "  +1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.17"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.597"
$ws.Range("E7").Value = "  +5.71%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.574"
$ws.Range("E9").Value = "  +12.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.54"
$ws.Range("E10").Value = "  +11.92%  "
$ws.Range("E11").Value = "  +6.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.24"
$ws.Range("E12").Value = "  +15.87%  "
$ws.Range("D13").Value = "2.981.20"
$ws.Range("E13").Value = "  +8.47%  "
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "2.595.83"
$ws.Range("E15").Value = "  +9.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.901"
$ws.Range("E16").Value = "  +9.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.83"
$ws.Range("E17").Value = "  +7.67%  "
$ws.Range("D18").Value = "46.285.36"
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("E19").Value = "  +6.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.97"
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("E21").Value = "  +9.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.07"
$ws.Range("E22").Value = "  +6.47%  "
$ws.Range("E23").Value = "  +4.13%  "
$ws.Range("E24").Value = "  +7.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  +13.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.07"
$ws.Range("E26").Value = "  +33.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +7.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.72"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.27"
$ws.Range("E30").Value = "  +2.48%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("E31").Value = "  +9.31%  "
$ws.Range("B32").Value = "LidoDAOToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.68"
$ws.Range("E32").Value = "  -1.91%  "
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.28"
$ws.Range("E34").Value = "  +18.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.76"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("E36").Value = "  +6.91%  "
$ws.Range("E37").Value = "  +3.66%  "
$ws.Range("E38").Value = "  +5.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.17"
$ws.Range("E39").Value = "  +8.71%  "
$ws.Range("E40").Value = "  +7.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  +12.76%  "
$ws.Range("E42").Value = "  +7.36%  "
$ws.Range("D43").Value = "2.057.75"
$ws.Range("E43").Value = "  +5.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.92"
$ws.Range("E44").Value = "  +40.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.81"
$ws.Range("E46").Value = "  -4.67%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.26"
$ws.Range("E47").Value = "  +9.20%  "
$ws.Range("E48").Value = "  -2.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.36"
$ws.Range("E49").Value = "  +9.21%  "
$ws.Range("E50").Value = "  +8.02%  "
$ws.Range("D51").Value = "2.844.11"
$ws.Range("E51").Value = "  +8.63%  "
